$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the weekly data block (row 104),
# pushing the existing rows 104-123 down to 106-125.
$ws.Range("A104:R105").Insert()

# Row 104: new weekly entry - Camote, 1a (guarda)
$ws.Cells.Item(104, 1).Value = 11
$ws.Cells.Item(104, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(104, 3).Value = "Bíobío"
$ws.Cells.Item(104, 4).Value = 44474
$ws.Cells.Item(104, 5).Value = 8
$ws.Cells.Item(104, 6).Value = 100112045
$ws.Cells.Item(104, 7).Value = "Zapallo"
$ws.Cells.Item(104, 8).Value = "Camote"
$ws.Cells.Item(104, 9).Value = "1a (guarda)"
$ws.Cells.Item(104, 10).Value = 600
$ws.Cells.Item(104, 11).Value = 600
$ws.Cells.Item(104, 12).Value = 650
$ws.Cells.Item(104, 13).Value = 625
$ws.Cells.Item(104, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(104, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(104, 16).Value = 625
$ws.Cells.Item(104, 17).Value = 1
$ws.Cells.Item(104, 18).Value = "Hortaliza"

# Row 105: new weekly entry - Camote, 2a (guarda)
$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44474
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = 100112045
$ws.Cells.Item(105, 7).Value = "Zapallo"
$ws.Cells.Item(105, 8).Value = "Camote"
$ws.Cells.Item(105, 9).Value = "2a (guarda)"
$ws.Cells.Item(105, 10).Value = 300
$ws.Cells.Item(105, 11).Value = 550
$ws.Cells.Item(105, 12).Value = 550
$ws.Cells.Item(105, 13).Value = 550
$ws.Cells.Item(105, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(105, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(105, 16).Value = 550
$ws.Cells.Item(105, 17).Value = 1
$ws.Cells.Item(105, 18).Value = "Hortaliza"
